# Regenerate save_data: K column (G) recomputed using "K instead of Strike#"
# Update the K values (column G) for rows 2-23 with the newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 2
    11 = 0
    12 = 2
    13 = 3
    14 = 1
    15 = 2
    16 = 0
    17 = 1
    18 = 1
    19 = 3
    20 = 1
    21 = 1
    22 = 1
    23 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

$wb.Save()
